# Update res_bus/vm_pu.xlsx values for the "case with 380 kV" run.
# Bus 0's per-unit voltage setpoint (column B) drops from 1.05 to 1.02,
# and the resulting voltage magnitudes recomputed for every bus column
# (B:F and I:N) on data rows 2-25. Columns A, G, H and row 1 (headers)
# are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.02, 1.075917653254934, 1.077441080185268, 1.079644051196198, 1.086508070661421, 1.056957280067858, 1.08081971452275, 1.080122372221687, 1.082319562943982, 1.089165734391528, 1.082354602489127)
    3 = @(1.02, 1.077476754437072, 1.078733602338909, 1.081155950515136, 1.088004396434375, 1.057507951031075, 1.08203561043378, 1.081231429478976, 1.083647885778299, 1.090479826711426, 1.083572225111815)
    4 = @(1.02, 1.078483749226014, 1.079568180000709, 1.08213273317307, 1.088971170208338, 1.057861897894416, 1.082820099061313, 1.081946738484381, 1.084505370813436, 1.091328165806895, 1.084357827803182)
    5 = @(1.02, 1.078906656432871, 1.079918619308929, 1.082543017431431, 1.089377262772752, 1.058010132330279, 1.083149359873499, 1.082246904057856, 1.084865379284689, 1.0916843439549, 1.084687556203487)
    6 = @(1.02, 1.078977639259271, 1.079977435344624, 1.082611885370039, 1.089445427839564, 1.058034988524886, 1.083204612854296, 1.082297271089718, 1.084925798462034, 1.091744120854767, 1.084742887649855)
    7 = @(1.02, 1.078489401824658, 1.079572864218638, 1.082138216797001, 1.088976597759986, 1.057863880825224, 1.08282450076655, 1.081950751466084, 1.084510183134214, 1.091332926889531, 1.084362235759345)
    8 = @(1.02, 1.076444946855246, 1.07787826478536, 1.080155323690313, 1.087014065694139, 1.057143875598344, 1.081231107995954, 1.080497668581685, 1.082768900397133, 1.089610249555251, 1.08276658018825)
    9 = @(1.02, 1.072827798667251, 1.074878299906563, 1.076649228934066, 1.083544385621846, 1.055856792826152, 1.078405590235939, 1.077919058013526, 1.079684666520258, 1.086559288018947, 1.079937049868924)
    10 = @(1.02, 1.070406011470404, 1.072868575061067, 1.074303265263095, 1.081223077230624, 1.054986181972485, 1.076509531590343, 1.076187419510009, 1.077617370279644, 1.084514516556759, 1.078038298602571)
    11 = @(1.02, 1.069354767889192, 1.071995934173745, 1.073285290922532, 1.080215872842643, 1.054606169556463, 1.075685481869607, 1.075434531018395, 1.076719460730531, 1.083626445350412, 1.077213078636711)
    12 = @(1.02, 1.068963888522817, 1.071671425068803, 1.07290683600602, 1.079841432847277, 1.054464556073674, 1.075378927196914, 1.075154404519396, 1.076385513848036, 1.08329616592812, 1.076906088621221)
    13 = @(1.02, 1.069047751783796, 1.071741050315712, 1.072988031141067, 1.079921766107145, 1.054494953531294, 1.075444705437615, 1.075214514007691, 1.076457165890621, 1.083367030692191, 1.076971960274568)
    14 = @(1.02, 1.06932246591689, 1.071969117775317, 1.073254014590629, 1.080184928087924, 1.054594473150266, 1.075660151519373, 1.07541138531181, 1.076691865265334, 1.083599152770442, 1.077187712314474)
    15 = @(1.02, 1.069491672934189, 1.072109588196414, 1.07341785122656, 1.080347028253255, 1.054655729391722, 1.075792832963605, 1.075532621679999, 1.076836414930411, 1.083742116217812, 1.077320582181576)
    16 = @(1.02, 1.070475723459542, 1.072926437734638, 1.074370778585936, 1.081289877798663, 1.055011337885748, 1.076564156192843, 1.076237320732543, 1.077676902721415, 1.084573397875508, 1.078093000778272)
    17 = @(1.02, 1.071092290190045, 1.073438172971076, 1.074967940385948, 1.081880743479881, 1.055233586804461, 1.077047165899908, 1.076678530373074, 1.078203374046666, 1.085094117392847, 1.078576696414522)
    18 = @(1.02, 1.071451673622714, 1.073736426860728, 1.075316047355702, 1.082225187081107, 1.055362928577604, 1.077328604192335, 1.076935584080573, 1.078510190510195, 1.08539758710192, 1.078858534381606)
    19 = @(1.02, 1.07157417209639, 1.073838084524619, 1.075434707955574, 1.082342600091611, 1.055406981339691, 1.077424517921941, 1.077023182710781, 1.078614762262127, 1.085501019075213, 1.07895458431971)
    20 = @(1.02, 1.071026164302686, 1.07338329274413, 1.074903892100985, 1.081817369770155, 1.055209771860639, 1.076995373908379, 1.076631223470022, 1.078146916168889, 1.085038275781617, 1.078524830872425)
    21 = @(1.02, 1.069241580601766, 1.071901967913325, 1.073175698392116, 1.080107442339063, 1.054565179835388, 1.07559672094926, 1.075353424662316, 1.076622763908077, 1.083530809999797, 1.077124191665676)
    22 = @(1.02, 1.068117218188398, 1.070968447068106, 1.072087177478201, 1.079030490118409, 1.054157235900805, 1.074714632356793, 1.074547297536566, 1.075662015240502, 1.082580628502966, 1.076240850406231)
    23 = @(1.02, 1.068713487947366, 1.071463531022956, 1.072664409941486, 1.079601581722378, 1.054373748632338, 1.075182503029835, 1.074974901637106, 1.076171561936503, 1.083084566064946, 1.076709385509294)
    24 = @(1.02, 1.071056044515934, 1.073408091475348, 1.074932833405081, 1.081846006234785, 1.055220533711803, 1.077018777385786, 1.076652600339244, 1.078172427873735, 1.085063508996879, 1.078548267585454)
    25 = @(1.02, 1.073764701459689, 1.07565554943525, 1.077557110863829, 1.084442788075469, 1.056191731076236, 1.079138203777018, 1.078587875751756, 1.080483943826566, 1.087349904284271, 1.080670703805263)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    $cols = @(2,3,4,5,6,9,10,11,12,13,14)
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item([int]$r, $cols[$i]).Value = $rowVals[$i]
    }
}
